# Fill empty ER sheets into every template (except Imaging)
# Adds a new "GEO_RNASEQ" worksheet (the ER / ontology-term lookup sheet)
# as the last tab of the workbook, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# --- add the new worksheet as the LAST tab -------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "GEO_RNASEQ"

# --- header row (row 1) ---------------------------------------------------
$ws.Cells.Item(1, 1).Value = ""
$ws.Cells.Item(1, 2).Value = "TermSourceRef"
$ws.Cells.Item(1, 3).Value = "Ontology"
$ws.Cells.Item(1, 4).Value = "TAN"
$ws.Cells.Item(1, 5).Value = "Content type (validation)"
$ws.Cells.Item(1, 6).Value = "Notes during templating"
$ws.Cells.Item(1, 7).Value = "Target term"
$ws.Cells.Item(1, 8).Value = "Instruction"
$ws.Cells.Item(1, 9).Value = "Requirement (m/o/n)"
$ws.Cells.Item(1, 10).Value = "Value (cv/s/d)"
$ws.Cells.Item(1, 11).Value = "Additional information"
$ws.Cells.Item(1, 12).Value = "Review comments"

# --- data rows 2-9 (column A = header names from the main sheet) ---------
$ws.Cells.Item(2, 1).Value = "Source Name"
$ws.Cells.Item(3, 1).Value = "Sample Name"

$ws.Cells.Item(4, 1).Value = "Parameter [Bio entity]"
$ws.Cells.Item(4, 2).Value = "NFDI4PSO:0000012"
$ws.Cells.Item(4, 3).Value = "NFDI4PSO"
$ws.Cells.Item(4, 4).Value = "http://purl.obolibrary.org/obo/NFDI4PSO_0000012"

$ws.Cells.Item(5, 1).Value = "Parameter [Biosource amount]"
$ws.Cells.Item(5, 2).Value = "NFDI4PSO:0000013"
$ws.Cells.Item(5, 3).Value = "NFDI4PSO"
$ws.Cells.Item(5, 4).Value = "http://purl.obolibrary.org/obo/NFDI4PSO_0000013"

$ws.Cells.Item(6, 1).Value = "Parameter [Extraction method]"
$ws.Cells.Item(6, 2).Value = "NFDI4PSO:0000054"
$ws.Cells.Item(6, 3).Value = "NFDI4PSO"
$ws.Cells.Item(6, 4).Value = "http://purl.obolibrary.org/obo/NFDI4PSO_0000054"

$ws.Cells.Item(7, 1).Value = "Parameter [Extraction buffer]"
$ws.Cells.Item(7, 2).Value = "NFDI4PSO:0000050"
$ws.Cells.Item(7, 3).Value = "NFDI4PSO"
$ws.Cells.Item(7, 4).Value = "http://purl.obolibrary.org/obo/NFDI4PSO_0000050"

$ws.Cells.Item(8, 1).Value = "Parameter [Extraction buffer volume]"
$ws.Cells.Item(8, 2).Value = "NFDI4PSO:0000051"
$ws.Cells.Item(8, 3).Value = "NFDI4PSO"
$ws.Cells.Item(8, 4).Value = "http://purl.obolibrary.org/obo/NFDI4PSO_0000051"

$ws.Cells.Item(9, 1).Value = "Parameter [RNA quality check]"
$ws.Cells.Item(9, 2).Value = "NFDI4PSO:0000062"
$ws.Cells.Item(9, 3).Value = "NFDI4PSO"
$ws.Cells.Item(9, 4).Value = "http://purl.obolibrary.org/obo/NFDI4PSO_0000062"

# --- column widths (best-effort match of the authored template) ----------
$ws.Columns.Item(1).ColumnWidth = 34.71
$ws.Columns.Item(2).ColumnWidth = 17.57
$ws.Columns.Item(3).ColumnWidth = 9.86
$ws.Columns.Item(4).ColumnWidth = 46.57
$ws.Columns.Item(5).ColumnWidth = 23.57
$ws.Columns.Item(6).ColumnWidth = 22.86
$ws.Columns.Item(7).ColumnWidth = 11.29
$ws.Columns.Item(8).ColumnWidth = 10.57
$ws.Columns.Item(9).ColumnWidth = 20.43
$ws.Columns.Item(10).ColumnWidth = 13.57
$ws.Columns.Item(11).ColumnWidth = 21.57
$ws.Columns.Item(12).ColumnWidth = 17.43

# --- selection / active-tab bookkeeping, matching a freshly-added sheet --
$ws.Cells.Select()
$ws.Activate()
